$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: add the missing "stop" date (C21) for the Draconid 2021 season ---
$ws.Range("B21").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("C21").Value = 44441

# --- Row 22: M3_09 Dryad 2021 (9.3) ---
$ws.Range("A17:F17").Copy()
$ws.Range("A22:F22").PasteSpecial(-4122)
$ws.Range("A22").Value = "M3_09 Dryad 2021"
$ws.Range("B22").Value = 44441
$ws.Range("C22").Value = 44474
$ws.Range("D22").Value = "'9.3"
$ws.Range("E22").Value = "Only balance changes/bug fixes"
$ws.Range("D17").Copy()
$ws.Range("D22").PasteSpecial(-4122)

# --- Row 23: M3_10 Cat 2021 (9.4, new expansion) ---
$ws.Range("A18:E18").Copy()
$ws.Range("A23:E23").PasteSpecial(-4122)
$ws.Range("A23").Value = "M3_10 Cat 2021"
$ws.Range("B23").Value = 44474
$ws.Range("C23").Value = 44497
$ws.Range("D23").Value = "'9.4"
$ws.Range("E23").Value = "New Expansion! Harvest of Sorrow"
$ws.Range("D18").Copy()
$ws.Range("D23").PasteSpecial(-4122)

# --- Row 24: M3_11 Mahakam 2021 (9.5, start Regis journey) ---
$ws.Range("A18:E18").Copy()
$ws.Range("A24:E24").PasteSpecial(-4122)
$ws.Range("A24").Value = "M3_11 Mahakam 2021"
$ws.Range("B24").Value = 44497
$ws.Range("C24").Value = 44537
$ws.Range("D24").Value = "'9.5"
$ws.Range("E24").Value = "Start Regis Journey"
$ws.Range("D18").Copy()
$ws.Range("D24").PasteSpecial(-4122)

# --- Row 25: M3_12 Wild Hunt 2021 (9.6, 12 new legendary cards) ---
$ws.Range("A18:E18").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)
$ws.Range("A25").Value = "M3_12 Wild Hunt 2021"
$ws.Range("B25").Value = 44537
$ws.Range("C25").Value = 44574
$ws.Range("D25").Value = "'9.6"
$ws.Range("E25").Value = "12 New Legendary Cards added"
$ws.Range("D18").Copy()
$ws.Range("D25").PasteSpecial(-4122)

# --- Row 26: M4_01 Wolf 2022 (10.1, draft out of early access) ---
$ws.Range("A18:E18").Copy()
$ws.Range("A26:E26").PasteSpecial(-4122)
$ws.Range("A26").Value = "M4_01 Wolf 2022"
$ws.Range("B26").Value = 44574
$ws.Range("C26").Value = 44600
$ws.Range("D26").Value = "'10.1"
$ws.Range("E26").Value = "Draft out of early access"
$ws.Range("D18").Copy()
$ws.Range("D26").PasteSpecial(-4122)

# --- Update the active selection / view to match what was left selected ---
$ws.Range("F29").Select()
